$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last checked" timestamp (column D) for all data rows (2:55)
$ws.Range("D2:D55").Value = 45974.307534722226

# Rebuild rows 19:55 (station, terminal, last-charge-end, and the D column already set above)
$data = New-Object 'object[,]' 37,3
$data[0,0] = '长沙特来电飞狐四方坪南区充电站'
$data[0,1] = '406号直流'
$data[0,2] = 45971.54614583333
$data[1,0] = '长沙特来电飞狐四方坪南区充电站'
$data[1,1] = '101号直流'
$data[1,2] = 45971.970555555556
$data[2,0] = '长沙特来电飞狐四方坪西区充电站'
$data[2,1] = '505号直流'
$data[2,2] = 45972.035127314812
$data[3,0] = '长沙特来电飞狐四方坪西区充电站'
$data[3,1] = '702号直流'
$data[3,2] = 45972.123761574076
$data[4,0] = '长沙特来电飞狐四方坪东区充电站'
$data[4,1] = '401号直流'
$data[4,2] = 45972.159884259258
$data[5,0] = '长沙特来电飞狐四方坪东区充电站'
$data[5,1] = '103号直流'
$data[5,2] = 45972.656111111108
$data[6,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[6,1] = '102号直流'
$data[6,2] = 45972.758530092593
$data[7,0] = '长沙特来电飞狐四方坪西区充电站'
$data[7,1] = '604号直流'
$data[7,2] = 45973.024872685186
$data[8,0] = '长沙特来电飞狐四方坪西区充电站'
$data[8,1] = '903号直流'
$data[8,2] = 45973.042685185188
$data[9,0] = '长沙特来电飞狐四方坪东区充电站'
$data[9,1] = '101号直流'
$data[9,2] = 45973.067175925928
$data[10,0] = '长沙特来电飞狐四方坪西区充电站'
$data[10,1] = '603号直流'
$data[10,2] = 45973.251481481479
$data[11,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[11,1] = '110号直流'
$data[11,2] = 45973.307268518518
$data[12,0] = '长沙特来电飞狐四方坪东区充电站'
$data[12,1] = '502号直流'
$data[12,2] = 45973.327638888892
$data[13,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[13,1] = '107号直流'
$data[13,2] = 45973.493518518517
$data[14,0] = '长沙特来电飞狐四方坪西区充电站'
$data[14,1] = '804号直流'
$data[14,2] = 45973.519178240742
$data[15,0] = '长沙特来电飞狐四方坪西区充电站'
$data[15,1] = '902号直流'
$data[15,2] = 45973.529386574075
$data[16,0] = '长沙特来电飞狐四方坪西区充电站'
$data[16,1] = '503号直流'
$data[16,2] = 45973.534062500003
$data[17,0] = '长沙特来电飞狐四方坪南区充电站'
$data[17,1] = '401号直流'
$data[17,2] = 45973.538935185185
$data[18,0] = '长沙特来电飞狐四方坪西区充电站'
$data[18,1] = '802号直流'
$data[18,2] = 45973.539178240739
$data[19,0] = '长沙特来电飞狐四方坪西区充电站'
$data[19,1] = '405号直流'
$data[19,2] = 45973.542581018519
$data[20,0] = '长沙特来电飞狐四方坪南区充电站'
$data[20,1] = '201号直流'
$data[20,2] = 45973.549907407411
$data[21,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[21,1] = '210号直流'
$data[21,2] = 45973.586909722224
$data[22,0] = '长沙特来电飞狐四方坪南区充电站'
$data[22,1] = '105号直流'
$data[22,2] = 45973.588495370372
$data[23,0] = '长沙特来电飞狐四方坪西区充电站'
$data[23,1] = 'A05号直流'
$data[23,2] = 45973.591238425928
$data[24,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[24,1] = '106号直流'
$data[24,2] = 45973.604259259257
$data[25,0] = '长沙特来电飞狐四方坪西区充电站'
$data[25,1] = '703号直流'
$data[25,2] = 45973.606365740743
$data[26,0] = '长沙特来电飞狐四方坪南区充电站'
$data[26,1] = '206号直流'
$data[26,2] = 45973.617754629631
$data[27,0] = '长沙特来电飞狐四方坪南区充电站'
$data[27,1] = '301号直流'
$data[27,2] = 45973.638402777775
$data[28,0] = '长沙特来电飞狐四方坪西区充电站'
$data[28,1] = 'B03号直流'
$data[28,2] = 45973.639351851853
$data[29,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[29,1] = '208号直流'
$data[29,2] = 45973.645648148151
$data[30,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[30,1] = '207号直流'
$data[30,2] = 45973.646122685182
$data[31,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[31,1] = '204号直流'
$data[31,2] = 45973.658449074072
$data[32,0] = '长沙特来电飞狐四方坪西区充电站'
$data[32,1] = '404号直流'
$data[32,2] = 45973.661736111113
$data[33,0] = '长沙特来电飞狐四方坪东区充电站'
$data[33,1] = '402号直流'
$data[33,2] = 45973.66201388889
$data[34,0] = '长沙特来电飞狐四方坪西区充电站'
$data[34,1] = '501号直流'
$data[34,2] = 45973.676249999997
$data[35,0] = '长沙特来电飞狐四方坪南区充电站'
$data[35,1] = '104号直流'
$data[35,2] = 45973.739710648151
$data[36,0] = '长沙特来电飞狐四方坪南区充电站'
$data[36,1] = '202号直流'
$data[36,2] = 45973.789178240739

$ws.Range("A19:C55").Value = $data

# Update the active cell selection
$ws.Range("F9").Select()
